# Update the GSMArena demo data: replace the three sample phones with newer
# Samsung Galaxy S10 family models, restyle the updated model cell, and move
# the active selection the way the author left it when they committed.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GSMArena")
$ws2 = $wb.Worksheets.Item("Calculator")

# --- Row 2: Galaxy S10+ -----------------------------------------------
$ws1.Range("D2").Value = "Galaxy S10+"
$ws1.Range("E2").Value = "6.4"""
$ws1.Range("F2").Value = "16MP"
$ws1.Range("G2").Value = "12GB RAM"
$ws1.Range("H2").Value = "4100mAh"

# --- Row 3: Galaxy S10 --------------------------------------------------
$ws1.Range("D3").Value = "Galaxy S10"
$ws1.Range("E3").Value = "6.1"""
$ws1.Range("F3").Value = "16MP"
$ws1.Range("G3").Value = "8GB RAM"
$ws1.Range("H3").Value = "3400mAh"

# --- Row 4: Galaxy View2 -------------------------------------------------
$ws1.Range("D4").Value = "Galaxy View2"
$ws1.Range("E4").Value = "17.3"""
$ws1.Range("F4").Value = "NO"
$ws1.Range("G4").Value = "3GB RAM"
$ws1.Range("H4").Value = "12000mAh"

# The new model name in D2 picks up a distinct monospace style (plain,
# no inherited number format/alignment - just the new font).
$ws1.Range("D2").ClearFormats() | Out-Null
$ws1.Range("D2").Font.Name = "Menlo"
$ws1.Range("D2").Font.Color = 2236962

# The workbook now opens on the GSMArena sheet with H5 selected, and the
# Calculator sheet is left on its previous E4 selection (no longer active).
$ws1.Activate() | Out-Null
$ws1.Range("H5").Select() | Out-Null

Write-Output "done"
